# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos
# table on the active sheet, row by row, matching the latest scrape.
# "" for the price means that row's price text did not change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  Price = "28.047.24";    Volume = "  -3.27%  " },
    @{ Row = 3;  Price = "1.860.74";     Volume = "  -2.77%  " },
    @{ Row = 4;  Price = "1.005";        Volume = "  +0.26%  " },
    @{ Row = 5;  Price = "318.10";       Volume = "  -2.06%  " },
    @{ Row = 6;  Price = "1.004";        Volume = "  +0.26%  " },
    @{ Row = 7;  Price = "0.4364";       Volume = "  -4.98%  " },
    @{ Row = 8;  Price = "0.3677";       Volume = "  -3.90%  " },
    @{ Row = 9;  Price = "0.07478";      Volume = "  -3.22%  " },
    @{ Row = 10; Price = "0.9332";       Volume = "  -4.97%  " },
    @{ Row = 11; Price = "21.31";        Volume = "  -3.79%  " },
    @{ Row = 12; Price = "1.893.68";     Volume = "  -0.75%  " },
    @{ Row = 13; Price = "6.683";        Volume = "  -3.85%  " },
    @{ Row = 14; Price = "5.421";        Volume = "  -4.39%  " },
    @{ Row = 15; Price = "0.06901";      Volume = "  -1.85%  " },
    @{ Row = 16; Price = "1.005";        Volume = "  +0.20%  " },
    @{ Row = 17; Price = "81.30";        Volume = "  -3.12%  " },
    @{ Row = 18; Price = "0.000008982";  Volume = "  -5.15%  " },
    @{ Row = 19; Price = "";             Volume = "  +0.31%  " },
    @{ Row = 20; Price = "15.80";        Volume = "  -5.57%  " },
    @{ Row = 21; Price = "28.038.58";    Volume = "  -3.20%  " },
    @{ Row = 22; Price = "5.104";        Volume = "  -4.11%  " },
    @{ Row = 23; Price = "10.79";        Volume = "  -1.17%  " },
    @{ Row = 24; Price = "2.111.31";     Volume = "  -1.26%  " },
    @{ Row = 25; Price = "2.012";        Volume = "  -3.67%  " },
    @{ Row = 26; Price = "154.25";       Volume = "  -2.70%  " },
    @{ Row = 27; Price = "18.31";        Volume = "  -4.05%  " },
    @{ Row = 28; Price = "5.348";        Volume = "  -5.68%  " },
    @{ Row = 29; Price = "112.94";       Volume = "  -3.94%  " },
    @{ Row = 30; Price = "1.725";        Volume = "  -7.08%  " },
    @{ Row = 31; Price = "0.08968";      Volume = "  -3.60%  " },
    @{ Row = 32; Price = "0.7924";       Volume = "  -8.70%  " },
    @{ Row = 33; Price = "4.810";        Volume = "  -5.34%  " },
    @{ Row = 34; Price = "3.044";        Volume = "  +0.63%  " },
    @{ Row = 35; Price = "1.168";        Volume = "  -6.85%  " },
    @{ Row = 36; Price = "";             Volume = "  +0.29%  " },
    @{ Row = 37; Price = "1.119";        Volume = "  -3.07%  " },
    @{ Row = 38; Price = "0.05401";      Volume = "  -5.77%  " },
    @{ Row = 39; Price = "0.01950";      Volume = "  -4.50%  " },
    @{ Row = 40; Price = "2.933";        Volume = "  +2.51%  " },
    @{ Row = 41; Price = "0.5223";       Volume = "  -5.25%  " },
    @{ Row = 42; Price = "6.984";        Volume = "  -5.71%  " },
    @{ Row = 43; Price = "0.1677";       Volume = "  -4.44%  " },
    @{ Row = 44; Price = "8.704";        Volume = "  -6.78%  " },
    @{ Row = 45; Price = "0.06714";      Volume = "  -2.31%  " },
    @{ Row = 46; Price = "0.4852";       Volume = "  -6.60%  " },
    @{ Row = 47; Price = "10.56";        Volume = "  -6.33%  " },
    @{ Row = 48; Price = "106.66";       Volume = "  -3.52%  " },
    @{ Row = 49; Price = "1.003";        Volume = "  +0.22%  " },
    @{ Row = 50; Price = "1.902";        Volume = "  -7.52%  " },
    @{ Row = 51; Price = "1.666";        Volume = "  -6.65%  " }
)

foreach ($u in $updates) {
    if ($u.Price -ne "") {
        $priceCell = $ws.Range("D" + $u.Row)
        # Prefix with an apostrophe so a price string that happens to look
        # numeric (e.g. "1.005") is kept as literal text instead of being
        # parsed into a float/date by the COM value-setter, matching the
        # inline-string storage used throughout this column.
        $priceCell.Value = "'" + $u.Price
        # Restore the cell's original (default/no-op) style now that the
        # text is committed, so the quote-prefix formatting applied above
        # doesn't linger on the cell.
        $priceCell.Style = "Normal"
    }
    $ws.Range("E" + $u.Row).Value = $u.Volume
}
